$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.800.62"
$ws.Range("E2").Value = "  -3.34%  "
$ws.Range("D3").Value = "2.273.50"
$ws.Range("E3").Value = "  -3.75%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "530.38"
$ws.Range("E5").Value = "  -5.04%  "
$ws.Range("D6").Value = "130.53"
$ws.Range("E6").Value = "  -2.09%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").Value = "0.582"
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("D9").Value = "2.272.52"
$ws.Range("E9").Value = "  -3.74%  "
$ws.Range("D10").Value = "0.0986"
$ws.Range("E10").Value = "  -5.99%  "
$ws.Range("E11").Value = "  -3.56%  "
$ws.Range("E12").Value = "  -0.34%  "
$ws.Range("D13").Value = "0.328"
$ws.Range("E13").Value = "  -3.81%  "
$ws.Range("D14").Value = "23.39"
$ws.Range("E14").Value = "  -3.43%  "
$ws.Range("D15").Value = "2.681.87"
$ws.Range("E15").Value = "  -3.78%  "
$ws.Range("D16").Value = "57.803.01"
$ws.Range("E16").Value = "  -3.28%  "
$ws.Range("E17").Value = "  -4.48%  "
$ws.Range("D18").Value = "2.314.63"
$ws.Range("E18").Value = "  -3.13%  "
$ws.Range("D19").Value = "10.48"
$ws.Range("E19").Value = "  -5.18%  "
$ws.Range("D20").Value = "4.20"
$ws.Range("E20").Value = "  -5.76%  "
$ws.Range("D21").Value = "310.84"
$ws.Range("E21").Value = "  -2.83%  "
$ws.Range("E22").Value = "  -4.38%  "
$ws.Range("E23").Value = "  -0.06%  "
$ws.Range("D24").Value = "62.47"
$ws.Range("E24").Value = "  -2.67%  "
$ws.Range("E25").Value = "  -2.82%  "
$ws.Range("E26").Value = "  +0.19%  "
$ws.Range("D27").Value = "7.95"
$ws.Range("E27").Value = "  -5.18%  "
$ws.Range("E28").Value = "  -7.09%  "
$ws.Range("D29").Value = "170.49"
$ws.Range("E29").Value = "  -0.03%  "
$ws.Range("E30").Value = "  -6.08%  "
$ws.Range("D31").Value = "0.0₃0714"
$ws.Range("E31").Value = "  -5.47%  "
$ws.Range("D32").Value = "5.72"
$ws.Range("E32").Value = "  -5.79%  "
$ws.Range("E33").Value = "  -6.13%  "
$ws.Range("D34").Value = "0.377"
$ws.Range("E34").Value = "  -5.05%  "
$ws.Range("E35").Value = "  +0.01%  "
$ws.Range("D36").Value = "17.68"
$ws.Range("E36").Value = "  -2.33%  "
$ws.Range("E37").Value = "  -0.06%  "
$ws.Range("E38").Value = "  -7.13%  "
$ws.Range("D39").Value = "3.87"
$ws.Range("E39").Value = "  -5.72%  "
$ws.Range("D40").Value = "38.19"
$ws.Range("E40").Value = "  -1.01%  "
$ws.Range("E41").Value = "  -6.47%  "
$ws.Range("D42").Value = "141.09"
$ws.Range("E42").Value = "  -2.41%  "
$ws.Range("D43").Value = "285.62"
$ws.Range("E43").Value = "  -9.98%  "
$ws.Range("E44").Value = "  -3.54%  "
$ws.Range("D45").Value = "0.0945"
$ws.Range("E45").Value = "  -2.01%  "
$ws.Range("E46").Value = "  -3.29%  "
$ws.Range("E47").Value = "  -3.34%  "
$ws.Range("D48").Value = "18.01"
$ws.Range("E48").Value = "  -6.88%  "
$ws.Range("E49").Value = "  -3.68%  "
$ws.Range("E50").Value = "  -1.25%  "
$ws.Range("E51").Value = "  -0.51%  "
